$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 38, shifting existing rows 38-41 down to 39-42
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new record's data
$ws.Range("A38").Value = 7
$ws.Range("B38").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C38").Value = "Ñuble"
$ws.Range("D38").Value = 45258
$ws.Range("E38").Value = 16
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100107
$ws.Range("H38").Value = "Otros"
$ws.Range("I38").Value = 100107002
$ws.Range("J38").Value = "Chirimoya"
$ws.Range("K38").Value = "Cultivar IV Región"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 60
$ws.Range("N38").Value = 20000
$ws.Range("O38").Value = 20000
$ws.Range("P38").Value = 20000
$ws.Range("Q38").Value = "$/bandeja 10 kilos"
$ws.Range("R38").Value = "Provincia de Limarí"
$ws.Range("S38").Value = 2000
$ws.Range("T38").Value = 10
